$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20 (shifts old rows 20-28 down to 21-29)
$ws.Rows.Item(20).Insert()

# Copy the style used by the other data rows (e.g. row 21, col A) onto the new A20 cell
$ws.Range("A21").Copy()
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new LIN.F row (VWAP / brand data merge)
$ws.Range("A20").Value = "LIN.F"
$ws.Range("B20").Value = 505
$ws.Range("C20").Value = 179.922
$ws.Range("D20").Value = 179.25
$ws.Range("E20").Value = 21.202
$ws.Range("F20").Value = 166.35
$ws.Range("G20").Value = 196.95
$ws.Range("H20").Value = 134.65
$ws.Range("I20").Value = 224.1
